$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.1074906481627241
$ws.Range("D2").Value = 0.9153740072163092

$ws.Range("C3").Value = -0.6264229408834185
$ws.Range("D3").Value = 0.5374838426647335

$ws.Range("C4").Value = -2.5650094385363
$ws.Range("D4").Value = 0.01765955831445232

$ws.Range("C5").Value = -1.140051697774375
$ws.Range("D5").Value = 0.2665250379685453

$ws.Range("C6").Value = -0.7939661827392596
$ws.Range("D6").Value = 0.4356960541109605

$ws.Range("C7").Value = -1.95594457460606
$ws.Range("D7").Value = 0.06328770853159993

$ws.Range("C8").Value = -1.419506315513625
$ws.Range("D8").Value = 0.1697664104934751

$ws.Range("C9").Value = -1.497717579089876
$ws.Range("D9").Value = 0.1484181780557379

$ws.Range("C10").Value = -1.130700962789817
$ws.Range("D10").Value = 0.270357360651881

$ws.Range("C11").Value = 0.4189060076883179
$ws.Range("D11").Value = 0.6793464733782986
